$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.801.93"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").Value = "'1.700.41"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("D4").Value = "'0.9966"
$ws.Range("E4").Value = "  -0.86%  "

$ws.Range("D5").Value = "'314.43"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").Value = "'0.9969"
$ws.Range("E6").Value = "  -0.66%  "

$ws.Range("D7").Value = "'0.3972"
$ws.Range("E7").Value = "  +1.92%  "

$ws.Range("D8").Value = "'0.4030"
$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").Value = "'1.471"
$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'52.60"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D11").Value = "'0.9952"
$ws.Range("E11").Value = "  -1.08%  "

$ws.Range("D12").Value = "'0.08839"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").Value = "'26.15"
$ws.Range("E13").Value = "  +5.90%  "

$ws.Range("D14").Value = "'7.512"
$ws.Range("E14").Value = "  -0.08%  "

$ws.Range("D15").Value = "'8.001"
$ws.Range("E15").Value = "  -0.28%  "

$ws.Range("D16").Value = "'0.00001350"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").Value = "'1.695.44"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "'95.98"
$ws.Range("E18").Value = "  -2.49%  "

$ws.Range("D19").Value = "'0.07167"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("D20").Value = "'20.81"
$ws.Range("E20").Value = "  +4.51%  "

$ws.Range("D21").Value = "'7.383"
$ws.Range("E21").Value = "  +1.43%  "

$ws.Range("D22").Value = "'0.9977"
$ws.Range("E22").Value = "  -0.59%  "

$ws.Range("D23").Value = "'14.40"
$ws.Range("E23").Value = "  +1.12%  "

$ws.Range("D24").Value = "'24.795.65"
$ws.Range("E24").Value = "  +1.57%  "

$ws.Range("D25").Value = "'2.364"
$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("D26").Value = "'2.933"
$ws.Range("E26").Value = "  -1.88%  "

$ws.Range("D27").Value = "'23.43"
$ws.Range("E27").Value = "  +3.45%  "

$ws.Range("D28").Value = "'6.368"
$ws.Range("E28").Value = "  +22.65%  "

$ws.Range("D29").Value = "'161.15"
$ws.Range("E29").Value = "  -2.10%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'145.51"
$ws.Range("E30").Value = "  +5.66%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'8.473"
$ws.Range("E31").Value = "  -1.65%  "

$ws.Range("D32").Value = "'2.402"
$ws.Range("E32").Value = "  +22.37%  "

$ws.Range("D33").Value = "'1.878.60"
$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("D34").Value = "'0.08619"
$ws.Range("E34").Value = "  -1.93%  "

$ws.Range("D35").Value = "'7.271"
$ws.Range("E35").Value = "  -2.81%  "

$ws.Range("D36").Value = "'0.03136"
$ws.Range("E36").Value = "  +6.28%  "

$ws.Range("D37").Value = "'1.037"
$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("D38").Value = "'0.2848"
$ws.Range("E38").Value = "  +1.43%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.8368"
$ws.Range("E39").Value = "  +5.44%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.09463"
$ws.Range("E40").Value = "  +3.48%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'10.75"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("D42").Value = "'14.12"
$ws.Range("E42").Value = "  -0.60%  "

$ws.Range("D43").Value = "'1.477"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("D44").Value = "'17.48"
$ws.Range("E44").Value = "  +1.04%  "

$ws.Range("D45").Value = "'2.692"
$ws.Range("E45").Value = "  +3.02%  "

$ws.Range("D46").Value = "'0.7442"
$ws.Range("E46").Value = "  +2.89%  "

$ws.Range("E47").Value = "  +0.16%  "

$ws.Range("D48").Value = "'1.374"
$ws.Range("E48").Value = "  -1.30%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.08628"
$ws.Range("E49").Value = "  +7.36%  "

$ws.Range("B50").Value = "Frax"
$ws.Range("C50").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D50").Value = "'0.9963"
$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("E51").Value = "  +0.51%  "
